# Add season-record columns (Wins / Losses / Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): AC1=Wins, AD1=Losses, AE1=Ties
# Copy the existing header formatting (bold font, border, centered) from AB1
# so the new header cells share the same style index as the rest of row 1.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows 2-43: constant season record for every player row
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 29).Value = 98
    $ws.Cells.Item($row, 30).Value = 64
    $ws.Cells.Item($row, 31).Value = 0
}
